# Adicionando BR2019, BR2018 e BR2017
# Appends two new "de-para" rows to Planilha1: CSA / CSA, and Fortaleza / FOR.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$ws.Range("A55").Value = "CSA"
$ws.Range("B55").Value = "CSA"

$ws.Range("A56").Value = "Fortaleza"
$ws.Range("B56").Value = "FOR"

# Match the author's final selection/view position on the new last row.
$ws.Range("B56").Select()
